{"js": "// Update CATCH.DEP.cod / CATCH.DEP.hake values in the \"Table 15. Socioeconomic\n// Factors\" table. Cells are addressed by (row, column) index so that the\n// update is unambiguous even though several of the new values collide with\n// other (old or new) values elsewhere in the table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// columns: 0 COUNTRIES, 1 FLEET.MOBILITY, 2 CATCH.DEP.cod, 3 CATCH.DEP.hake, 4 ADAPTIVE.MNG\nconst COL_COD = 2;\nconst COL_HAKE = 3;\n\n// rows are keyed by country -> table row index (row 0 is the header row)\nconst rowByCountry = {\n  BE: 1, DK: 2, DE: 3, EE: 4, IE: 5, ES: 6, FR: 7,\n  LV: 8, LT: 9, NL: 10, PL: 11, PT: 12, FI: 13, SE: 14\n};\n\n// [country, column, oldValue, newValue]\nconst changes = [\n  [\"BE\", COL_COD, \"0.959\", \"0.970\"],\n  [\"BE\", COL_HAKE, \"0.954\", \"0.996\"],\n  [\"DK\", COL_COD, \"0.908\", \"0.959\"],\n  [\"DK\", COL_HAKE, \"0.897\", \"0.995\"],\n  [\"DE\", COL_COD, \"0.840\", \"0.764\"],\n  [\"DE\", COL_HAKE, \"0.989\", \"0.999\"],\n  [\"EE\", COL_COD, \"0.734\", \"0.948\"],\n  [\"IE\", COL_COD, \"0.904\", \"0.989\"],\n  [\"ES\", COL_COD, \"0.604\", \"0.871\"],\n  [\"ES\", COL_HAKE, \"0.825\", \"0.916\"],\n  [\"FR\", COL_COD, \"0.939\", \"0.977\"],\n  [\"FR\", COL_HAKE, \"0.827\", \"0.838\"],\n  [\"LV\", COL_COD, \"0.697\", \"0.714\"],\n  [\"LT\", COL_COD, \"0.692\", \"0.464\"],\n  [\"NL\", COL_COD, \"0.996\", \"1.000\"],\n  [\"PL\", COL_COD, \"0.554\", \"0.401\"],\n  [\"PT\", COL_COD, \"0.650\", \"0.839\"],\n  [\"PT\", COL_HAKE, \"0.773\", \"0.939\"],\n  [\"FI\", COL_COD, \"0.687\", \"0.921\"],\n  [\"SE\", COL_COD, \"0.647\", \"0.722\"],\n  [\"SE\", COL_HAKE, \"0.567\", \"0.996\"],\n];\n\n// Locate, for each change, the text range holding the old value inside the\n// correct cell (scoping the search to the cell keeps every lookup unique\n// even after earlier replacements have introduced new matching digits\n// elsewhere in the table).\nconst searchResults = [];\nfor (const [country, col, oldValue] of changes) {\n  const row = rowByCountry[country];\n  const cell = table.getCell(row, col);\n  const found = cell.body.search(oldValue, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  searchResults.push(found);\n}\nawait context.sync();\n\nfor (let i = 0; i < changes.length; i++) {\n  const [, , , newValue] = changes[i];\n  const found = searchResults[i];\n  if (found.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for change index \" + i + \", got \" + found.items.length);\n  }\n  found.items[0].insertText(newValue, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update CATCH.DEP.cod / CATCH.DEP.hake values in the \"Table 15. Socioeconomic\n# Factors\" table. Cells are addressed by their (row, column) index - rather\n# than via Find/Replace - so the update is unambiguous even though several\n# of the new values collide with other (old or new) values elsewhere in the\n# table, and Find.Execute on a Range in this environment is not reliably\n# constrained to that Range.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Columns (1-based, as used by Table.Cell): 1 COUNTRIES, 2 FLEET.MOBILITY,\n# 3 CATCH.DEP.cod, 4 CATCH.DEP.hake, 5 ADAPTIVE.MNG\n$COL_COD = 3\n$COL_HAKE = 4\n\n# Rows (1-based; row 1 is the header row)\n$rowByCountry = @{\n  BE = 2; DK = 3; DE = 4; EE = 5; IE = 6; ES = 7; FR = 8;\n  LV = 9; LT = 10; NL = 11; PL = 12; PT = 13; FI = 14; SE = 15\n}\n\n# country, column, old value (for sanity-check), new value\n$changes = @(\n  @(\"BE\", $COL_COD,  \"0.959\", \"0.970\"),\n  @(\"BE\", $COL_HAKE, \"0.954\", \"0.996\"),\n  @(\"DK\", $COL_COD,  \"0.908\", \"0.959\"),\n  @(\"DK\", $COL_HAKE, \"0.897\", \"0.995\"),\n  @(\"DE\", $COL_COD,  \"0.840\", \"0.764\"),\n  @(\"DE\", $COL_HAKE, \"0.989\", \"0.999\"),\n  @(\"EE\", $COL_COD,  \"0.734\", \"0.948\"),\n  @(\"IE\", $COL_COD,  \"0.904\", \"0.989\"),\n  @(\"ES\", $COL_COD,  \"0.604\", \"0.871\"),\n  @(\"ES\", $COL_HAKE, \"0.825\", \"0.916\"),\n  @(\"FR\", $COL_COD,  \"0.939\", \"0.977\"),\n  @(\"FR\", $COL_HAKE, \"0.827\", \"0.838\"),\n  @(\"LV\", $COL_COD,  \"0.697\", \"0.714\"),\n  @(\"LT\", $COL_COD,  \"0.692\", \"0.464\"),\n  @(\"NL\", $COL_COD,  \"0.996\", \"1.000\"),\n  @(\"PL\", $COL_COD,  \"0.554\", \"0.401\"),\n  @(\"PT\", $COL_COD,  \"0.650\", \"0.839\"),\n  @(\"PT\", $COL_HAKE, \"0.773\", \"0.939\"),\n  @(\"FI\", $COL_COD,  \"0.687\", \"0.921\"),\n  @(\"SE\", $COL_COD,  \"0.647\", \"0.722\"),\n  @(\"SE\", $COL_HAKE, \"0.567\", \"0.996\")\n)\n\nforeach ($change in $changes) {\n  $country = $change[0]\n  $col = $change[1]\n  $oldValue = $change[2]\n  $newValue = $change[3]\n  $row = $rowByCountry[$country]\n\n  $cell = $t.Cell($row, $col)\n  $r = $cell.Range\n  # Trim the trailing end-of-cell marker(s) so only the value text is compared/replaced.\n  [void]$r.MoveEnd(1, -1)\n  if ($r.Text -ne $oldValue) {\n    throw \"Unexpected value in $country/$col : expected '$oldValue' got '$($r.Text)'\"\n  }\n  $r.Text = $newValue\n}\n"}
